$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename the existing three sheets and add a fourth, matching the new
#    LinkML-generated workbook layout:
#      NamedThing -> SampleCollection
#      Person -> Sample
#      PersonCollection -> Air Sample
#      (new) -> Soil Sample
# ---------------------------------------------------------------------------
$wsCollection = $wb.Worksheets.Item(1)
$wsSample     = $wb.Worksheets.Item(2)
$wsAir        = $wb.Worksheets.Item(3)

$wsCollection.Name = "SampleCollection"
$wsSample.Name     = "Sample"
$wsAir.Name         = "Air Sample"

# Copy (rather than Add) the "Air Sample" sheet to seed the new "Soil
# Sample" sheet so it inherits the same sheet formatting (margins,
# outline/page-setup props, default row height, ...) as the other
# LinkML-generated sheets instead of Excel's blank-sheet defaults.
$wsAir.Copy($null, $wsAir)
$wsSoil = $wb.Worksheets.Item(4)
$wsSoil.Name = "Soil Sample"

# ---------------------------------------------------------------------------
# 2. SampleCollection sheet: id, samples (was id, name, description)
# ---------------------------------------------------------------------------
$wsCollection.Range("A1").Value = "id"
$wsCollection.Range("B1").Value = "samples"
$wsCollection.Range("C1").ClearContents()

# ---------------------------------------------------------------------------
# 3. Sample sheet: id, latitude, longitude, species, sample biome
#    (was primary_email, birth_date, age_in_years, vital_status, id, name,
#     description)
# ---------------------------------------------------------------------------
$wsSample.Range("A1").Value = "id"
$wsSample.Range("B1").Value = "latitude"
$wsSample.Range("C1").Value = "longitude"
$wsSample.Range("D1").Value = "species"
$wsSample.Range("E1").Value = "sample biome"
$wsSample.Range("F1:G1").ClearContents()

# replace the old vital_status validation (D) with an empty species list,
# and add the new sample-biome list on E
$wsSample.Range("D2:D1048576").Validation.Delete()
$wsSample.Range("D2:D1048576").Validation.Add(3, 1, 1, '""')
$wsSample.Range("D2:D1048576").Validation.ShowInput = $false
$wsSample.Range("D2:D1048576").Validation.ShowError = $false

$wsSample.Range("E2:E1048576").Validation.Add(3, 1, 1, '"forest,lake,ocean,desert,air"')
$wsSample.Range("E2:E1048576").Validation.ShowInput = $false
$wsSample.Range("E2:E1048576").Validation.ShowError = $false

# ---------------------------------------------------------------------------
# 4. Air Sample sheet: altitude, id, latitude, longitude, species,
#    sample biome (was just "entries")
# ---------------------------------------------------------------------------
$wsAir.Range("A1").Value = "altitude"
$wsAir.Range("B1").Value = "id"
$wsAir.Range("C1").Value = "latitude"
$wsAir.Range("D1").Value = "longitude"
$wsAir.Range("E1").Value = "species"
$wsAir.Range("F1").Value = "sample biome"

$wsAir.Range("E2:E1048576").Validation.Add(3, 1, 1, '""')
$wsAir.Range("E2:E1048576").Validation.ShowInput = $false
$wsAir.Range("E2:E1048576").Validation.ShowError = $false

$wsAir.Range("F2:F1048576").Validation.Add(3, 1, 1, '"forest,lake,ocean,desert,air"')
$wsAir.Range("F2:F1048576").Validation.ShowInput = $false
$wsAir.Range("F2:F1048576").Validation.ShowError = $false

# ---------------------------------------------------------------------------
# 5. Soil Sample sheet (brand new): depth, id, latitude, longitude, species,
#    sample biome
# ---------------------------------------------------------------------------
$wsSoil.Range("A1").Value = "depth"
$wsSoil.Range("B1").Value = "id"
$wsSoil.Range("C1").Value = "latitude"
$wsSoil.Range("D1").Value = "longitude"
$wsSoil.Range("E1").Value = "species"
$wsSoil.Range("F1").Value = "sample biome"

$wsSoil.Range("E2:E1048576").Validation.Add(3, 1, 1, '""')
$wsSoil.Range("E2:E1048576").Validation.ShowInput = $false
$wsSoil.Range("E2:E1048576").Validation.ShowError = $false

$wsSoil.Range("F2:F1048576").Validation.Add(3, 1, 1, '"forest,lake,ocean,desert,air"')
$wsSoil.Range("F2:F1048576").Validation.ShowInput = $false
$wsSoil.Range("F2:F1048576").Validation.ShowError = $false

# ---------------------------------------------------------------------------
# Leave the selection/active sheet on the first sheet, matching the original
# workbook's default view.
# ---------------------------------------------------------------------------
$wsCollection.Activate()
$wsCollection.Range("A1").Select()
